$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 912.125
$ws.Range("I9").Value = 799.75
$ws.Range("K9").Value = 799.75
$ws.Range("M9").Value = -630.75
$ws.Range("H12").Value = 256.72726
$ws.Range("I12").Value = 192.2
$ws.Range("K12").Value = 192.2
$ws.Range("M12").Value = -22.19999999999999
$ws.Range("H33").Value = 340.14285
$ws.Range("I33").Value = 382.33334
$ws.Range("K33").Value = 382.33334
$ws.Range("M33").Value = -153.33334
$ws.Range("H40").Value = 18190.348
$ws.Range("I40").Value = 15898.8125
$ws.Range("J40").Value = 23428.143
$ws.Range("K40").Value = 15898.8125
$ws.Range("L40").Value = 23428.143
$ws.Range("M40").Value = -15723.8125
$ws.Range("N40").Value = -23778.143
$ws.Range("H55").Value = 530
$ws.Range("J55").Value = 353.75
$ws.Range("L55").Value = 353.75
$ws.Range("N55").Value = -781.75
$ws.Range("H106").Value = 1448.7727
$ws.Range("I106").Value = 1298.579
$ws.Range("K106").Value = 1298.579
$ws.Range("M106").Value = -667.579
$ws.Range("H125").Value = 2346.3572
$ws.Range("I125").Value = 1238.7778
$ws.Range("K125").Value = 11149.0002
$ws.Range("M125").Value = -8689.0002
$ws.Range("H132").Value = 167279.78
$ws.Range("I132").Value = 255190.58
$ws.Range("K132").Value = 765571.74
$ws.Range("M132").Value = -763041.74
$ws.Range("H135").Value = 4354.048
$ws.Range("I135").Value = 2046.2307
$ws.Range("K135").Value = 18416.0763
$ws.Range("M135").Value = -15881.0763
$ws.Range("H137").Value = 168959.58
$ws.Range("I137").Value = 184019.67
$ws.Range("J137").Value = 3298.6
$ws.Range("K137").Value = 552059.01
$ws.Range("L137").Value = 9895.799999999999
$ws.Range("M137").Value = -549509.01
$ws.Range("N137").Value = -14995.8
$ws.Range("H138").Value = 5456.485
$ws.Range("I138").Value = 1402.3125
$ws.Range("J138").Value = 6753.82
$ws.Range("K138").Value = 4206.9375
$ws.Range("L138").Value = 20261.46
$ws.Range("M138").Value = 933.0625
$ws.Range("N138").Value = -30541.46
$ws.Range("H141").Value = 2731.6938
$ws.Range("I141").Value = 2473.3489
$ws.Range("K141").Value = 7420.0467
$ws.Range("M141").Value = -2240.0467

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 16302.258
$ws.Range("I32").Value = 17153.06
$ws.Range("J32").Value = 13095.385
$ws.Range("K32").Value = 17153.06
$ws.Range("L32").Value = 13095.385
$ws.Range("M32").Value = -16866.06
$ws.Range("N32").Value = -13669.385
$ws.Range("H46").Value = 28128.889
$ws.Range("J46").Value = 28128.889
$ws.Range("L46").Value = 28128.889
$ws.Range("N46").Value = -28766.889
$ws.Range("H74").Value = 1830.875
$ws.Range("I74").Value = 1093.6923
$ws.Range("J74").Value = 3199.9285
$ws.Range("K74").Value = 1093.6923
$ws.Range("L74").Value = 3199.9285
$ws.Range("M74").Value = -219.6922999999999
$ws.Range("N74").Value = -4947.9285
$ws.Range("H77").Value = 1830.875
$ws.Range("I77").Value = 1093.6923
$ws.Range("J77").Value = 3199.9285
$ws.Range("K77").Value = 5468.461499999999
$ws.Range("L77").Value = 15999.6425
$ws.Range("M77").Value = -1100.461499999999
$ws.Range("N77").Value = -24735.6425
$ws.Range("H102").Value = 1315.7858
$ws.Range("I102").Value = 1385.3846
$ws.Range("K102").Value = 1385.3846
$ws.Range("M102").Value = 236.6153999999999
$ws.Range("H122").Value = 3871.413
$ws.Range("I122").Value = 2693.842
$ws.Range("J122").Value = 9464.875
$ws.Range("K122").Value = 8081.526
$ws.Range("L122").Value = 28394.625
$ws.Range("M122").Value = -5631.526
$ws.Range("N122").Value = -33294.625
$ws.Range("H132").Value = 14845.413
$ws.Range("I132").Value = 18491.697
$ws.Range("K132").Value = 55475.091
$ws.Range("M132").Value = -52945.091

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H26").Value = 13122.5
$ws.Range("I26").Value = 9643.200000000001
$ws.Range("K26").Value = 9643.200000000001
$ws.Range("M26").Value = -9351.200000000001
$ws.Range("H94").Value = 864.1667
$ws.Range("I94").Value = 664
$ws.Range("J94").Value = 1264.5
$ws.Range("K94").Value = 664
$ws.Range("L94").Value = 1264.5
$ws.Range("M94").Value = -213
$ws.Range("N94").Value = -2166.5
$ws.Range("H134").Value = 2618.5173
$ws.Range("I134").Value = 1405.5714
$ws.Range("K134").Value = 4216.7142
$ws.Range("M134").Value = -1681.7142

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 639
$ws.Range("J22").Value = 193.33333
$ws.Range("L22").Value = 193.33333
$ws.Range("N22").Value = -893.3333299999999
$ws.Range("H31").Value = 25645760
$ws.Range("J31").Value = 6487.778
$ws.Range("L31").Value = 6487.778
$ws.Range("N31").Value = -7077.778
$ws.Range("H34").Value = 25645760
$ws.Range("J34").Value = 6487.778
$ws.Range("L34").Value = 6487.778
$ws.Range("N34").Value = -6891.778
$ws.Range("H58").Value = 501929
$ws.Range("I58").Value = 1695.3077
$ws.Range("J58").Value = 1430934.4
$ws.Range("K58").Value = 1695.3077
$ws.Range("L58").Value = 1430934.4
$ws.Range("M58").Value = -1492.3077
$ws.Range("N58").Value = -1431340.4
$ws.Range("H105").Value = 2668.9333
$ws.Range("I105").Value = 3472.6667
$ws.Range("K105").Value = 3472.6667
$ws.Range("M105").Value = -1725.6667
$ws.Range("H109").Value = 37834.5
$ws.Range("J109").Value = 37834.5
$ws.Range("L109").Value = 37834.5
$ws.Range("N109").Value = -39914.5
$ws.Range("H134").Value = 2433.5264
$ws.Range("I134").Value = 2458.9375
$ws.Range("K134").Value = 7376.8125
$ws.Range("M134").Value = -4841.8125
$ws.Range("H136").Value = 501929
$ws.Range("I136").Value = 1695.3077
$ws.Range("J136").Value = 1430934.4
$ws.Range("K136").Value = 5085.9231
$ws.Range("L136").Value = 4292803.199999999
$ws.Range("M136").Value = -2535.9231
$ws.Range("N136").Value = -4297903.199999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 3871125.5
$ws.Range("I4").Value = 1143125.9
$ws.Range("J4").Value = 8645125
$ws.Range("K4").Value = 3429377.7
$ws.Range("L4").Value = 25935375
$ws.Range("M4").Value = -3429265.7
$ws.Range("N4").Value = -25935599
$ws.Range("H5").Value = 1238.3334
$ws.Range("I5").Value = 920.75
$ws.Range("J5").Value = 1601.2858
$ws.Range("K5").Value = 2762.25
$ws.Range("L5").Value = 4803.857400000001
$ws.Range("M5").Value = -2650.25
$ws.Range("N5").Value = -5027.857400000001
$ws.Range("H40").Value = 168.84616
$ws.Range("J40").Value = 217.75
$ws.Range("L40").Value = 871
$ws.Range("N40").Value = -1009
$ws.Range("H46").Value = 3697.5715
$ws.Range("J46").Value = 4996.4
$ws.Range("L46").Value = 14989.2
$ws.Range("N46").Value = -15171.2
$ws.Range("H135").Value = 1238.3334
$ws.Range("I135").Value = 920.75
$ws.Range("J135").Value = 1601.2858
$ws.Range("K135").Value = 8286.75
$ws.Range("L135").Value = 14411.5722
$ws.Range("M135").Value = -5751.75
$ws.Range("N135").Value = -19481.5722

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 101854.09
$ws.Range("I80").Value = 115377.89
$ws.Range("J80").Value = 40997
$ws.Range("K80").Value = 115377.89
$ws.Range("L80").Value = 40997
$ws.Range("M80").Value = -114379.89
$ws.Range("N80").Value = -42993
$ws.Range("H83").Value = 101854.09
$ws.Range("I83").Value = 115377.89
$ws.Range("J83").Value = 40997
$ws.Range("K83").Value = 576889.45
$ws.Range("L83").Value = 204985
$ws.Range("M83").Value = -571897.45
$ws.Range("N83").Value = -214969
$ws.Range("H123").Value = 31361.625
$ws.Range("J123").Value = 31361.625
$ws.Range("L123").Value = 31361.625
$ws.Range("N123").Value = -36261.625
$ws.Range("H132").Value = 378588.78
$ws.Range("I132").Value = 74708.64
$ws.Range("J132").Value = 2505749.8
$ws.Range("K132").Value = 224125.92
$ws.Range("L132").Value = 7517249.399999999
$ws.Range("M132").Value = -221595.92
$ws.Range("N132").Value = -7522309.399999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 5320.6
$ws.Range("I82").Value = 3101
$ws.Range("J82").Value = 6800.3335
$ws.Range("K82").Value = 3101
$ws.Range("L82").Value = 6800.3335
$ws.Range("M82").Value = -2740
$ws.Range("N82").Value = -7522.3335
$ws.Range("H85").Value = 5320.6
$ws.Range("I85").Value = 3101
$ws.Range("J85").Value = 6800.3335
$ws.Range("K85").Value = 3101
$ws.Range("L85").Value = 6800.3335
$ws.Range("M85").Value = -1853
$ws.Range("N85").Value = -9296.333500000001
$ws.Range("H132").Value = 2814.86
$ws.Range("I132").Value = 2324.7976
$ws.Range("K132").Value = 6974.3928
$ws.Range("M132").Value = -4444.3928
$ws.Range("H136").Value = 3114.5576
$ws.Range("I136").Value = 2136.1316
$ws.Range("K136").Value = 6408.3948
$ws.Range("M136").Value = -3858.3948
$ws.Range("H139").Value = 79299
$ws.Range("J139").Value = 79299
$ws.Range("L139").Value = 79299
$ws.Range("N139").Value = -89579

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 2991.6667
$ws.Range("I81").Value = 866.25
$ws.Range("J81").Value = 19995
$ws.Range("K81").Value = 1732.5
$ws.Range("L81").Value = 39990
$ws.Range("M81").Value = -671.5
$ws.Range("N81").Value = -42112
$ws.Range("H84").Value = 2991.6667
$ws.Range("I84").Value = 866.25
$ws.Range("J84").Value = 19995
$ws.Range("K84").Value = 8662.5
$ws.Range("L84").Value = 199950
$ws.Range("M84").Value = -3358.5
$ws.Range("N84").Value = -210558
$ws.Range("H96").Value = 2500
$ws.Range("I96").Value = 2500
$ws.Range("K96").Value = 2500
$ws.Range("M96").Value = -1127
$ws.Range("H136").Value = 8254.694
$ws.Range("I136").Value = 1932.5416
$ws.Range("J136").Value = 10742.099
$ws.Range("K136").Value = 5797.6248
$ws.Range("L136").Value = 32226.297
$ws.Range("M136").Value = -3247.6248
$ws.Range("N136").Value = -37326.297
